$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are safe (won't be auto-parsed as numbers/dates)
$ws.Range("D2").Value = '63.160.72'
$ws.Range("E2").Value = '  -4.42%  '
$ws.Range("D3").Value = '3.106.73'
$ws.Range("E3").Value = '  -4.59%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("E6").Value = '  -9.25%  '
$ws.Range("D8").Value = '3.092.10'
$ws.Range("E8").Value = '  -5.06%  '
$ws.Range("E9").Value = '  -4.54%  '
$ws.Range("E10").Value = '  -7.35%  '
$ws.Range("E11").Value = '  -9.34%  '
$ws.Range("E12").Value = '  -5.83%  '
$ws.Range("E13").Value = '  -8.31%  '
$ws.Range("E14").Value = '  -9.58%  '
$ws.Range("D15").Value = '3.616.63'
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = '63.298.79'
$ws.Range("E17").Value = '  -4.32%  '
$ws.Range("D18").Value = '3.108.35'
$ws.Range("E18").Value = '  -4.75%  '
$ws.Range("E19").Value = '  -7.94%  '
$ws.Range("E20").Value = '  -5.92%  '
$ws.Range("E21").Value = '  -5.54%  '
$ws.Range("E22").Value = '  -6.62%  '
$ws.Range("E23").Value = '  -4.98%  '
$ws.Range("E24").Value = '  -8.50%  '
$ws.Range("E25").Value = '  -3.76%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -8.99%  '
$ws.Range("E28").Value = '  -9.31%  '
$ws.Range("E29").Value = '  -12.68%  '
$ws.Range("E30").Value = '  -4.18%  '
$ws.Range("E31").Value = '  -14.56%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  -7.22%  '
$ws.Range("E34").Value = '  -6.63%  '
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  -8.50%  '
$ws.Range("E37").Value = '  -5.43%  '
$ws.Range("D38").Value = '0.0₃0736'
$ws.Range("E38").Value = '  -6.46%  '
$ws.Range("E39").Value = '  -9.11%  '
$ws.Range("E40").Value = '  -16.94%  '
$ws.Range("E41").Value = '  -8.28%  '
$ws.Range("E42").Value = '  -9.63%  '
$ws.Range("E43").Value = '  -5.94%  '
$ws.Range("D44").Value = '2.826.17'
$ws.Range("E45").Value = '  -9.85%  '
$ws.Range("E46").Value = '  -12.91%  '
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("E49").Value = '  -10.64%  '
$ws.Range("E50").Value = '  -5.90%  '
$ws.Range("E51").Value = '  -2.58%  '

# Cells whose new values look numeric: force text format first so Excel
# stores them as strings (matching the original inline-string cell type)
# instead of auto-converting to a Number.
$textCells = @("D5", "D6", "D9", "D11", "D20", "D21", "D22", "D25", "D28", "D29", "D31", "D34", "D36", "D37", "D39", "D40", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '604.00'
$ws.Range("D6").Value = '142.79'
$ws.Range("D9").Value = '0.515'
$ws.Range("D11").Value = '5.17'
$ws.Range("D20").Value = '469.47'
$ws.Range("D21").Value = '14.44'
$ws.Range("D22").Value = '0.698'
$ws.Range("D25").Value = '83.18'
$ws.Range("D28").Value = '8.30'
$ws.Range("D29").Value = '2.06'
$ws.Range("D31").Value = '0.113'
$ws.Range("D34").Value = '25.94'
$ws.Range("D36").Value = '5.86'
$ws.Range("D37").Value = '52.55'
$ws.Range("D39").Value = '450.80'
$ws.Range("D40").Value = '2.88'
$ws.Range("D45").Value = '0.262'
$ws.Range("D47").Value = '2.39'
$ws.Range("D48").Value = '0.999'
$ws.Range("D49").Value = '25.68'
$ws.Range("D51").Value = '118.04'
